# "Fruta / hortaliza, semanal"
#
# The weekly refresh re-shuffles which source record (Fecha, Volumen,
# Precio minimo/maximo/promedio, Origen, Variedad, Calidad, etc.) lands on
# which sheet row. Every data row (2-21) keeps the same fixed columns
# (Mercado ID, Mercado, Region, Codreg, Categoria ID, Categoria, Unidad de
# comercializacion, Kg o Unidades, Clasificacion) but the full per-record
# tuple of the remaining columns moves to a different row.
#
# Snapshot every row's full A:R tuple first (so sources aren't clobbered
# before they're read), then write each tuple to its new destination row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 21

# Capture the "before" state of every data row as a whole record.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapshot[$r] = $ws.Range("A$r`:R$r").Value()
}

# destinationRow = sourceRow : content that used to live at sourceRow now
# belongs at destinationRow.
$rowMap = @{
    2  = 18
    3  = 10
    4  = 20
    5  = 12
    6  = 13
    7  = 8
    8  = 5
    9  = 7
    10 = 6
    11 = 21
    12 = 3
    13 = 11
    14 = 19
    15 = 2
    16 = 15
    17 = 9
    18 = 16
    19 = 14
    20 = 17
    21 = 4
}

foreach ($destRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$destRow]
    $ws.Range("A$destRow`:R$destRow").Value = $snapshot[$sourceRow]
}
